$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = "<siie>"
$ws.Range("C3").Value = 49

# Row 4
$ws.Range("B4").Value = "<she>"
$ws.Range("C4").Value = 54

# Row 5
$ws.Range("B5").Value = "<on>"
$ws.Range("C5").Value = 50

# Row 6
$ws.Range("B6").Value = "<in>"
$ws.Range("C6").Value = 51

# Row 7
$ws.Range("B7").Value = "<rich>"
$ws.Range("C7").Value = 53

# Row 8
$ws.Range("C8").Value = 41

# Row 9
$ws.Range("B9").Value = "<its>"
